$d = $word.ActiveDocument

$d.Content.Find.Execute("002", $true, $false, $false, $false, $false, $true, 1, $false, "008", 2)
$d.Content.Find.Execute("Rakesh Singh", $true, $false, $false, $false, $false, $true, 1, $false, "Anirban Dutta", 2)
$d.Content.Find.Execute("Ajay Singh", $true, $false, $false, $false, $false, $true, 1, $false, "Sagar Dutta", 2)
$d.Content.Find.Execute("12347N", $true, $false, $false, $false, $false, $true, 1, $false, "1254A", 2)
$d.Content.Find.Execute("005", $true, $false, $false, $false, $false, $true, 1, $false, "12", 2)
$d.Content.Find.Execute("A+", $true, $false, $false, $false, $false, $true, 1, $false, "A", 2)
$d.Content.Find.Execute("09/04/2021", $true, $false, $false, $false, $false, $true, 1, $false, "12/02/2020", 2)
